$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 12:22"

# Update country data rows (values refreshed + a few countries swapped/reordered
# due to re-sorting by case totals)

$ws.Cells.Item(31, 1).Value = "Rumania"
$ws.Cells.Item(31, 2).Value = 8746
$ws.Cells.Item(31, 3).Value = 328
$ws.Cells.Item(31, 4).Value = 1892
$ws.Cells.Item(31, 5).Value = 6420
$ws.Cells.Item(31, 6).Value = 256
$ws.Cells.Item(31, 7).Value = 13
$ws.Cells.Item(31, 8).Value = 434

$ws.Cells.Item(50, 1).Value = "Finlandia"
$ws.Cells.Item(50, 2).Value = 3783
$ws.Cells.Item(50, 3).Value = 102
$ws.Cells.Item(50, 4).Value = 1700
$ws.Cells.Item(50, 5).Value = 1989
$ws.Cells.Item(50, 6).Value = 68
$ws.Cells.Item(50, 7).Value = 4
$ws.Cells.Item(50, 8).Value = 94

$ws.Cells.Item(56, 1).Value = "Marruecos"
$ws.Cells.Item(56, 2).Value = 2820
$ws.Cells.Item(56, 3).Value = 135
$ws.Cells.Item(56, 4).Value = 322
$ws.Cells.Item(56, 5).Value = 2360
$ws.Cells.Item(56, 6).Value = 1
$ws.Cells.Item(56, 7).Value = 1
$ws.Cells.Item(56, 8).Value = 138

$ws.Cells.Item(57, 1).Value = "Tailandia"
$ws.Cells.Item(57, 2).Value = 2765
$ws.Cells.Item(57, 3).Value = 32
$ws.Cells.Item(57, 4).Value = 1928
$ws.Cells.Item(57, 5).Value = 790
$ws.Cells.Item(57, 6).Value = 61
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 47

$ws.Cells.Item(116, 1).Value = "Vietnam"
$ws.Cells.Item(116, 2).Value = 268
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 203
$ws.Cells.Item(116, 5).Value = 65
$ws.Cells.Item(116, 6).Value = 8
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 0

$ws.Cells.Item(119, 1).Value = "Sri Lanka"
$ws.Cells.Item(119, 2).Value = 256
$ws.Cells.Item(119, 3).Value = 2
$ws.Cells.Item(119, 4).Value = 91
$ws.Cells.Item(119, 5).Value = 158
$ws.Cells.Item(119, 6).Value = 1
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 7

$ws.Cells.Item(139, 1).Value = "Etiopia"
$ws.Cells.Item(139, 2).Value = 108
$ws.Cells.Item(139, 3).Value = 3
$ws.Cells.Item(139, 4).Value = 16
$ws.Cells.Item(139, 5).Value = 89
$ws.Cells.Item(139, 6).Value = 1
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 3

$ws.Cells.Item(140, 1).Value = "Birmania"
$ws.Cells.Item(140, 2).Value = 107
$ws.Cells.Item(140, 3).Value = 9
$ws.Cells.Item(140, 4).Value = 5
$ws.Cells.Item(140, 5).Value = 97
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 5

$ws.Cells.Item(154, 1).Value = "Zambia"
$ws.Cells.Item(154, 2).Value = 61
$ws.Cells.Item(154, 3).Value = 4
$ws.Cells.Item(154, 4).Value = 33
$ws.Cells.Item(154, 5).Value = 26
$ws.Cells.Item(154, 6).Value = 1
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 2

$ws.Cells.Item(155, 1).Value = "Cabo Verde"
$ws.Cells.Item(155, 2).Value = 58
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 1
$ws.Cells.Item(155, 5).Value = 56
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 1

$ws.Cells.Item(159, 1).Value = "Maldivas"
$ws.Cells.Item(159, 2).Value = 51
$ws.Cells.Item(159, 3).Value = 16
$ws.Cells.Item(159, 4).Value = 16
$ws.Cells.Item(159, 5).Value = 35
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

$ws.Cells.Item(160, 1).Value = "Libia"
$ws.Cells.Item(160, 2).Value = 49
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 11
$ws.Cells.Item(160, 5).Value = 37
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 1

$ws.Cells.Item(161, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(161, 2).Value = 46
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 0
$ws.Cells.Item(161, 5).Value = 46
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 0

$ws.Cells.Item(162, 1).Value = "Macao"
$ws.Cells.Item(162, 2).Value = 45
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 16
$ws.Cells.Item(162, 5).Value = 29
$ws.Cells.Item(162, 6).Value = 1
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 0

$ws.Cells.Item(163, 1).Value = "Haiti"
$ws.Cells.Item(163, 2).Value = 44
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 0
$ws.Cells.Item(163, 5).Value = 41
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 3

$ws.Cells.Item(164, 1).Value = "Puerto Rico"
$ws.Cells.Item(164, 2).Value = 39
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 1
$ws.Cells.Item(164, 5).Value = 36
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 2

$ws.Cells.Item(165, 1).Value = "Eritrea"
$ws.Cells.Item(165, 2).Value = 39
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 3
$ws.Cells.Item(165, 5).Value = 36
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 0

$ws.Cells.Item(166, 1).Value = "Siria"
$ws.Cells.Item(166, 2).Value = 38
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 5
$ws.Cells.Item(166, 5).Value = 31
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 2

$ws.Cells.Item(167, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(167, 2).Value = 37
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 19
$ws.Cells.Item(167, 5).Value = 16
$ws.Cells.Item(167, 6).Value = 5
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 2

$ws.Cells.Item(184, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(184, 2).Value = 17
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 17
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(185, 1).Value = "Fiyi"
$ws.Cells.Item(185, 2).Value = 17
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 17
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Sudan del Sur"
$ws.Cells.Item(210, 2).Value = 4
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 0
$ws.Cells.Item(210, 5).Value = 4
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(211, 2).Value = 4
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 0
$ws.Cells.Item(211, 5).Value = 4
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0
